$wb = $excel.ActiveWorkbook

# --- Sheet "Input Data": update vendor-name sample values ---
$wsInput = $wb.Worksheets.Item("Input Data")
$wsInput.Range("A2").Value = "Vendor Satu"
$wsInput.Range("A3").Value = "Vendor Dua"

# Re-apply center horizontal alignment across the sample data rows
$wsInput.Range("A2:J3").HorizontalAlignment = -4108   # xlCenter

# Update the remembered selection on this sheet
$wsInput.Range("I11").Select()

# --- Sheet "Keterangan": re-apply center alignment to the category list ---
$wsKet = $wb.Worksheets.Item("Keterangan")
$wsKet.Range("B3:C18").HorizontalAlignment = -4108   # xlCenter
